# Set column L ("Diferencia Stock") equal to column K ("Stock Mínimo Objetivo")
# for each data row (3 through 50), and update the "Total_Ajuste_Stock" summary
# cell (C64) to the sum of column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

for ($row = 3; $row -le 50; $row++) {
    $kValue = $ws.Cells.Item($row, 11).Value2  # column K = 11
    $ws.Cells.Item($row, 12).Value2 = $kValue  # column L = 12
}

# Update the summary metric "Total_Ajuste_Stock" to match the new total of column L
$ws.Range("C64").Value2 = $ws.Application.WorksheetFunction.Sum($ws.Range("L3:L50"))
